$wb = $excel.ActiveWorkbook

# --- Metadata sheet: insert a new "Experiment description" column after column A ---
$wsMeta = $wb.Worksheets.Item("Metadata")
[void]$wsMeta.Columns.Item(2).Insert()
$wsMeta.Range("B1").Value = "Experiment description"
$wsMeta.Range("B2").Value = "Three concentration gradient for measuring v max and Km."

# --- Restore / update the selection on every other sheet (drop the stray
#     "K6:M8" range from the sqref, keep just the real active cell) ---
$wsGroups = $wb.Worksheets.Item("Groups")
[void]$wsGroups.Range("B7").Select()

$wsSpecies = $wb.Worksheets.Item("Species")
[void]$wsSpecies.Range("J9").Select()

$wsBaseMedia = $wb.Worksheets.Item("Base Media")
[void]$wsBaseMedia.Range("D11").Select()

$wsCarbonSource = $wb.Worksheets.Item("Carbon Source")
[void]$wsCarbonSource.Range("J13").Select()

$wsCSConcentration = $wb.Worksheets.Item("CS Concentration")
[void]$wsCSConcentration.Range("K6").Select()

$wsInhibitor = $wb.Worksheets.Item("Inhibitor")
[void]$wsInhibitor.Range("B2").Select()

$wsInhibitorConc = $wb.Worksheets.Item("Inhibitor Conc")
[void]$wsInhibitorConc.Range("J18").Select()

$wsComments = $wb.Worksheets.Item("Comments")
[void]$wsComments.Range("K2").Select()

# --- Make Metadata the active sheet/tab again, with B3 selected ---
[void]$wsMeta.Activate()
[void]$wsMeta.Range("B3").Select()
